# "Cambio la Fig. de Barras y reubico la Cita del Abstract"
# The bars-figure data table loses its last series (impl5) and every
# remaining series's TotalCellArea/TotalDynamicPower/CellLeakagePower
# numbers shift up one row; impl4's row then picks up the plain
# "Normal" formatting that impl5's row used to have. The selection /
# cursor position is also moved (the abstract citation was relocated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- shift the numeric series up one row -------------------------------
# row2 (impl1) now shows what used to be row3's (impl2's) numbers, etc.
$ws.Range("B2").Value = 357688.253875
$ws.Range("C2").Value = 699.9652
$ws.Range("D2").Value = 1.9611

$ws.Range("B3").Value = 290179.448967
$ws.Range("C3").Value = 661.2084
$ws.Range("D3").Value = 1.6008

$ws.Range("B4").Value = 260340.888442
$ws.Range("C4").Value = 660.2362
$ws.Range("D4").Value = 1.3871

# --- row 5 (impl4) takes on the plain/unstyled look that row 6 (impl5)
#     used to have, instead of the highlighted table style ---------------
$ws.Range("A1:D1").Copy() | Out-Null
$ws.Range("A5:D5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- the impl5 row (row 6) is dropped entirely --------------------------
$ws.Rows.Item(6).Delete() | Out-Null

# --- the columns get very slightly narrower now that the table shrank --
$ws.Columns.Item(1).ColumnWidth = 9.833333333333334
$ws.Columns.Item(2).ColumnWidth = 11.5
$ws.Columns.Item(3).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 15.833333333333334

$farRange1 = $ws.Range($ws.Cells.Item(1, 5), $ws.Cells.Item(1, 997))
$farRange1.EntireColumn.ColumnWidth = 7.666666666666667

$farRange2 = $ws.Range($ws.Cells.Item(1, 998), $ws.Cells.Item(1, 1025))
$farRange2.EntireColumn.ColumnWidth = 9.666666666666666

# --- reposition the selection (abstract citation moved) ----------------
$ws.Range("D17").Select() | Out-Null

Write-Host "powerArea figure/table updated"
